$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sku / price_net / pid rows appended after the existing data (row 57).
# Column A = sku, Column B = price_net, Column C = pid.
$newRows = @(
    @{ Row = 58; Sku = "FD9749-400";    Price = 549; Pid = "NI115O04G-K11"; Wrap = $true  },
    @{ Row = 59; Sku = "GSB550CA";      Price = 339; Pid = "NE216D05Q-A13"; Wrap = $false },
    @{ Row = 60; Sku = "FB1843-141";    Price = 449; Pid = "NI116D0II-T11"; Wrap = $true  },
    @{ Row = 61; Sku = "DD1391-300";    Price = 499; Pid = "NI112O0GN-M11"; Wrap = $true  },
    @{ Row = 62; Sku = "DR8057-500";    Price = 589; Pid = "JOC11A032-M11"; Wrap = $true  },
    @{ Row = 63; Sku = "DR9512-001";    Price = 549; Pid = "NI112N03R-Q11"; Wrap = $true  }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.Sku
    $ws.Range("B$rowNum").Value = $r.Price
    $ws.Range("C$rowNum").Value = $r.Pid

    if ($r.Wrap) {
        $ws.Range("A$rowNum").WrapText = $true
        $ws.Rows.Item($rowNum).RowHeight = 16
    }
}

# Update the view state to match the post-edit selection/scroll position.
[void]$excel.Goto($ws.Range("A26"), $true)
[void]$ws.Range("A64").Select()
